$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme
$c = $cs.Colors(1)
$c.RGB = $c.RGB
